$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels: remove the inner space from each compound word.
$ws.Range("B1").Value = "TOTALKGS "
$ws.Range("D1").Value = "TOTALDAYS "
$ws.Range("F1").Value = "TOTALAMOUNT"
$ws.Range("H1").Value = "NETPAY"

# Update the view: scroll back to top-left and move the selection to H1.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H1").Select()
